$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scaling")
$tbl = $ws.ListObjects.Item("tbl_I_scaling")
$newRow = $tbl.ListRows.Add()
$ws.Range("A9").Value = 2005
$ws.Range("B9").Value = "*"
$ws.Range("C9").Value = "Cd"
$ws.Range("D9").Value = "NL"
$ws.Range("E9").Value = "*"
$ws.Range("G9").Value = 0.8
